$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.311.70"
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").Value = "1.787.29"
$ws.Range("E3").Value = "  -2.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "340.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3454"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.194"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07455"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.461"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "

$ws.Range("D15").Value = "1.787.86"
$ws.Range("E15").Value = "  -2.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.117"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001091"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06695"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.19%  "

$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.498"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("D23").Value = "27.299.59"
$ws.Range("E23").Value = "  -0.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.389"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.470"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.496"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "156.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.02%  "

$ws.Range("D30").Value = "1.987.88"
$ws.Range("E30").Value = "  -2.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.033"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.973"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08816"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.42%  "

$ws.Range("E36").Value = "  -4.41%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02413"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.405"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6833"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06425"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2202"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.250"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.437"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6396"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.872"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.135"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07134"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.78%  "
